$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update quantities on the "Daily " requisition sheet:
#  - C11(4/64) Blue (row 11): Qty 50 -> 10
#  - C25Y(4/64) Blue (row 17): Qty (blank) -> 10
$ws.Range("D11").Value = 10
$ws.Range("D17").Value = 10

# Move the cursor / selection to the cell the author ended on
$ws.Range("L21").Select()
